# Merge the three separate runs "<id>", "p001v_1", "</id>" into a single
# run containing "<id>p001v_1</id>", keeping the formatting of the first
# run (Courier New, color 7f6000, sz/szCs 18). Find.Execute collapses a
# multi-run match into one run using the first run's character formatting
# when the replacement text is supplied, which is exactly the shape the
# diff expects.
$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p001v_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p001v_1</id>", 2)
